# Updated cryptos list on Fri Oct  6 13:53:11 UTC 2023 with GitHub Actions
#
# Refresh the Price (column D) and Volume(1h) (column E) figures for each
# coin row, and fix up two rows whose coins were re-ordered (ImmutableX /
# ARBITRUM and mCoin / MXToken swapped places).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rows whose Coin name / Link / Price / Volume all changed because two
# entries swapped positions in the source feed.
$swaps = @(
    @{ Row = 37; Coin = "ImmutableX"; Link = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"; Price = "0.558"; Volume = "  -1.64%  " },
    @{ Row = 38; Coin = "ARBITRUM";   Link = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb";   Price = "0.875"; Volume = "  -0.49%  " },
    @{ Row = 44; Coin = "mCoin";      Link = "https://coinranking.com/coin/fzVgyjBcRc9+mcoin-mcoin";   Price = "2.46";  Volume = "  +0.15%  " },
    @{ Row = 45; Coin = "MXToken";    Link = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx";   Price = "2.21";  Volume = "  -1.07%  " }
)

# Rows that only need their Price (column D) and Volume (column E) refreshed.
$priceUpdates = @(
    @{ Row = 2;  Price = "27.501.21"; Volume = "  -1.58%  " },
    @{ Row = 3;  Price = "1.629.25";  Volume = "  -0.66%  " },
    @{ Row = 5;  Price = "211.43";    Volume = "  -0.47%  " },
    @{ Row = 7;  Price = "1.00";      Volume = "  +0.10%  " },
    @{ Row = 8;  Price = "23.00";     Volume = "  -1.62%  " },
    @{ Row = 9;  Price = "0.262";     Volume = "  +0.12%  " },
    @{ Row = 12; Price = "1.861.76";  Volume = "  -0.56%  " },
    @{ Row = 13; Price = "1.632.95";  Volume = "  -0.44%  " },
    @{ Row = 16; Price = "64.97";     Volume = "  +0.59%  " },
    @{ Row = 17; Price = "27.523.21"; Volume = "  -1.35%  " },
    @{ Row = 18; Price = "228.66";    Volume = "  -2.06%  " },
    @{ Row = 22; Price = "10.73";     Volume = "  +7.26%  " },
    @{ Row = 23; Price = "4.36";      Volume = "  +1.09%  " },
    @{ Row = 25; Price = "149.04";    Volume = "  -0.99%  " },
    @{ Row = 26; Price = "6.86";      Volume = "  -1.27%  " },
    @{ Row = 28; Price = "15.57";     Volume = "  -0.70%  " },
    @{ Row = 31; Price = "0.0481";    Volume = "  -0.50%  " },
    @{ Row = 32; Price = "3.27";      Volume = "  -1.00%  " },
    @{ Row = 33; Price = "1.463.91";  Volume = "  -0.55%  " },
    @{ Row = 40; Price = "0.914";     Volume = "  -1.43%  " },
    @{ Row = 43; Price = "67.94";     Volume = "  -1.66%  " },
    @{ Row = 46; Price = "5.36";      Volume = "  -1.16%  " },
    @{ Row = 47; Price = "1.771.57";  Volume = "  -0.65%  " },
    @{ Row = 49; Price = "87.33";     Volume = "  +0.03%  " },
    @{ Row = 50; Price = "0.0993";    Volume = "  +0.17%  " },
    @{ Row = 51; Price = "0.0₇0986"; Volume = "  -7.14%  " }
)

# Rows that only need their Volume (column E) refreshed.
$volumeOnlyUpdates = @(
    @{ Row = 4;  Volume = "  +0.16%  " },
    @{ Row = 6;  Volume = "  -1.03%  " },
    @{ Row = 10; Volume = "  -0.27%  " },
    @{ Row = 11; Volume = "  -3.38%  " },
    @{ Row = 14; Volume = "  -0.36%  " },
    @{ Row = 15; Volume = "  -0.47%  " },
    @{ Row = 19; Volume = "  -0.64%  " },
    @{ Row = 20; Volume = "  -0.65%  " },
    @{ Row = 21; Volume = "  +0.10%  " },
    @{ Row = 24; Volume = "  +2.30%  " },
    @{ Row = 27; Volume = "  -1.10%  " },
    @{ Row = 29; Volume = "  +0.08%  " },
    @{ Row = 30; Volume = "  -0.76%  " },
    @{ Row = 34; Volume = "  -0.79%  " },
    @{ Row = 35; Volume = "  -1.14%  " },
    @{ Row = 36; Volume = "  -1.38%  " },
    @{ Row = 39; Volume = "  -0.39%  " },
    @{ Row = 41; Volume = "  +0.09%  " },
    @{ Row = 42; Volume = "  +0.62%  " },
    @{ Row = 48; Volume = "  +1.70%  " }
)

# Column D ("Price") contains strings that look numeric (e.g. "1.00",
# "23.00", "27.501.21"); setting .Value directly would let Excel coerce
# them into real numbers and lose the intended text formatting. Forcing
# the number format to Text before assigning - then restoring the cell's
# style - keeps the value as plain text without leaving a stray style
# behind.
function Set-TextValue($cell, [string]$value) {
    $cell.NumberFormat = "@"
    $cell.Value = $value
    $cell.Style = "Normal"
}

foreach ($item in $swaps) {
    $r = $item.Row
    $ws.Range("B$r").Value = $item.Coin
    $ws.Range("C$r").Value = $item.Link
    Set-TextValue $ws.Range("D$r") $item.Price
    $ws.Range("E$r").Value = $item.Volume
}

foreach ($item in $priceUpdates) {
    $r = $item.Row
    Set-TextValue $ws.Range("D$r") $item.Price
    $ws.Range("E$r").Value = $item.Volume
}

foreach ($item in $volumeOnlyUpdates) {
    $r = $item.Row
    $ws.Range("E$r").Value = $item.Volume
}
